$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1128.881
$ws.Range("I15").Value = 1128.881
$ws.Range("K15").Value = 3386.643
$ws.Range("M15").Value = -3217.643
$ws.Range("H29").Value = 1000
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 3000
$ws.Range("N29").Value = -3562
$ws.Range("H38").Value = 105.333336
$ws.Range("I38").Value = 105.333336
$ws.Range("K38").Value = 316.000008
$ws.Range("M38").Value = 55.99999200000002
$ws.Range("H58").Value = 2466.75
$ws.Range("J58").Value = 7000
$ws.Range("L58").Value = 21000
$ws.Range("N58").Value = -21300
$ws.Range("H63").Value = 50999.5
$ws.Range("J63").Value = 50999.5
$ws.Range("L63").Value = 50999.5
$ws.Range("N63").Value = -52247.5
$ws.Range("H66").Value = 50999.5
$ws.Range("J66").Value = 50999.5
$ws.Range("L66").Value = 152998.5
$ws.Range("N66").Value = -159238.5
$ws.Range("H70").Value = 1277008.5
$ws.Range("I70").Value = 10204081
$ws.Range("J70").Value = 1712.4286
$ws.Range("K70").Value = 30612243
$ws.Range("L70").Value = 5137.2858
$ws.Range("M70").Value = -30611973
$ws.Range("N70").Value = -5677.2858
$ws.Range("H73").Value = 1277008.5
$ws.Range("I73").Value = 10204081
$ws.Range("J73").Value = 1712.4286
$ws.Range("K73").Value = 30612243
$ws.Range("L73").Value = 5137.2858
$ws.Range("M73").Value = -30611307
$ws.Range("N73").Value = -7009.2858
$ws.Range("H76").Value = 4574980.5
$ws.Range("I76").Value = 67574.625
$ws.Range("K76").Value = 67574.625
$ws.Range("M76").Value = -67259.625
$ws.Range("H79").Value = 4574980.5
$ws.Range("I79").Value = 67574.625
$ws.Range("K79").Value = 67574.625
$ws.Range("M79").Value = -66482.625
$ws.Range("H80").Value = 949792.2
$ws.Range("I80").Value = 1625373.9
$ws.Range("J80").Value = 3977.7
$ws.Range("K80").Value = 4876121.699999999
$ws.Range("L80").Value = 11933.1
$ws.Range("M80").Value = -4875123.699999999
$ws.Range("N80").Value = -13929.1
$ws.Range("H83").Value = 949792.2
$ws.Range("I83").Value = 1625373.9
$ws.Range("J83").Value = 3977.7
$ws.Range("K83").Value = 14628365.1
$ws.Range("L83").Value = 35799.3
$ws.Range("M83").Value = -14623373.1
$ws.Range("N83").Value = -45783.3
$ws.Range("H87").Value = 121400
$ws.Range("J87").Value = 121400
$ws.Range("L87").Value = 121400
$ws.Range("N87").Value = -123896
$ws.Range("H88").Value = 2162.7646
$ws.Range("I88").Value = 1597.4
$ws.Range("J88").Value = 2398.3333
$ws.Range("K88").Value = 1597.4
$ws.Range("L88").Value = 2398.3333
$ws.Range("M88").Value = -1191.4
$ws.Range("N88").Value = -3210.3333
$ws.Range("H90").Value = 121400
$ws.Range("J90").Value = 121400
$ws.Range("L90").Value = 364200
$ws.Range("N90").Value = -376680
$ws.Range("H91").Value = 2162.7646
$ws.Range("I91").Value = 1597.4
$ws.Range("J91").Value = 2398.3333
$ws.Range("K91").Value = 1597.4
$ws.Range("L91").Value = 2398.3333
$ws.Range("M91").Value = -193.4000000000001
$ws.Range("N91").Value = -5206.3333
$ws.Range("H98").Value = 628.7027
$ws.Range("I98").Value = 630.75
$ws.Range("J98").Value = 555
$ws.Range("K98").Value = 630.75
$ws.Range("L98").Value = 555
$ws.Range("M98").Value = 867.25
$ws.Range("N98").Value = -3551
$ws.Range("H100").Value = 1473.52
$ws.Range("I100").Value = 1252.8823
$ws.Range("J100").Value = 1942.375
$ws.Range("K100").Value = 1252.8823
$ws.Range("L100").Value = 1942.375
$ws.Range("M100").Value = -711.8823
$ws.Range("N100").Value = -3024.375
$ws.Range("H112").Value = 2095.0334
$ws.Range("I112").Value = 2924.75
$ws.Range("K112").Value = 8774.25
$ws.Range("M112").Value = -7666.25
$ws.Range("H122").Value = 628.7027
$ws.Range("I122").Value = 630.75
$ws.Range("J122").Value = 555
$ws.Range("K122").Value = 1892.25
$ws.Range("L122").Value = 1665
$ws.Range("M122").Value = 557.75
$ws.Range("N122").Value = -6565
$ws.Range("H132").Value = 2756.1128
$ws.Range("I132").Value = 2578.984
$ws.Range("K132").Value = 7736.951999999999
$ws.Range("M132").Value = -5206.951999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2636.14
$ws.Range("J32").Value = 1599.5714
$ws.Range("L32").Value = 1599.5714
$ws.Range("N32").Value = -2173.5714
$ws.Range("H45").Value = 2174.25
$ws.Range("I45").Value = 1848.5
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 1848.5
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -1471.5
$ws.Range("N45").Value = -3254
$ws.Range("H61").Value = 1477.1428
$ws.Range("I61").Value = 1398.4615
$ws.Range("K61").Value = 1398.4615
$ws.Range("M61").Value = -1186.4615
$ws.Range("H74").Value = 3904.138
$ws.Range("I74").Value = 4031.04
$ws.Range("K74").Value = 4031.04
$ws.Range("M74").Value = -3157.04
$ws.Range("H77").Value = 3904.138
$ws.Range("I77").Value = 4031.04
$ws.Range("K77").Value = 20155.2
$ws.Range("M77").Value = -15787.2
$ws.Range("H97").Value = 683.8095
$ws.Range("I97").Value = 683.4211
$ws.Range("J97").Value = 687.5
$ws.Range("K97").Value = 683.4211
$ws.Range("L97").Value = 687.5
$ws.Range("M97").Value = -187.4211
$ws.Range("N97").Value = -1679.5
$ws.Range("H136").Value = 1477.1428
$ws.Range("I136").Value = 1398.4615
$ws.Range("K136").Value = 4195.3845
$ws.Range("M136").Value = -1645.3845

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3390.9614
$ws.Range("J94").Value = 3032.3333
$ws.Range("L94").Value = 3032.3333
$ws.Range("N94").Value = -3934.3333
$ws.Range("H105").Value = 1684.3158
$ws.Range("I105").Value = 1512.625
$ws.Range("K105").Value = 1512.625
$ws.Range("M105").Value = 234.375
$ws.Range("H134").Value = 1165.75
$ws.Range("I134").Value = 1165.75
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3497.25
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -962.25
$ws.Range("N134").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2405.25
$ws.Range("I99").Value = 2199
$ws.Range("J99").Value = 2749
$ws.Range("K99").Value = 2199
$ws.Range("L99").Value = 2749
$ws.Range("M99").Value = -701
$ws.Range("N99").Value = -5745
$ws.Range("H126").Value = 2405.25
$ws.Range("I126").Value = 2199
$ws.Range("J126").Value = 2749
$ws.Range("K126").Value = 6597
$ws.Range("L126").Value = 8247
$ws.Range("M126").Value = -4127
$ws.Range("N126").Value = -13187

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1060.1052
$ws.Range("J5").Value = 1423.2
$ws.Range("L5").Value = 4269.6
$ws.Range("N5").Value = -4493.6
$ws.Range("H23").Value = 61.333332
$ws.Range("I23").Value = 44.8
$ws.Range("J23").Value = 69.59999999999999
$ws.Range("K23").Value = 134.4
$ws.Range("L23").Value = 208.8
$ws.Range("M23").Value = 100.6
$ws.Range("N23").Value = -678.8
$ws.Range("H135").Value = 1060.1052
$ws.Range("J135").Value = 1423.2
$ws.Range("L135").Value = 12808.8
$ws.Range("N135").Value = -17878.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1541.0834
$ws.Range("I97").Value = 1588.5
$ws.Range("J97").Value = 1304
$ws.Range("K97").Value = 1588.5
$ws.Range("L97").Value = 1304
$ws.Range("M97").Value = -1092.5
$ws.Range("N97").Value = -2296
$ws.Range("H122").Value = 1619.1111
$ws.Range("I122").Value = 1474.6316
$ws.Range("J122").Value = 1962.25
$ws.Range("K122").Value = 4423.8948
$ws.Range("L122").Value = 5886.75
$ws.Range("M122").Value = -1973.8948
$ws.Range("N122").Value = -10786.75
$ws.Range("H132").Value = 6584.5
$ws.Range("J132").Value = 8709.429
$ws.Range("L132").Value = 26128.287
$ws.Range("N132").Value = -31188.287

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2578.3845
$ws.Range("I40").Value = 2002.1578
$ws.Range("J40").Value = 4142.4287
$ws.Range("K40").Value = 2002.1578
$ws.Range("L40").Value = 4142.4287
$ws.Range("M40").Value = -1866.1578
$ws.Range("N40").Value = -4414.4287
$ws.Range("H68").Value = 5357.08
$ws.Range("I68").Value = 5146.6875
$ws.Range("J68").Value = 5731.1113
$ws.Range("K68").Value = 5146.6875
$ws.Range("L68").Value = 5731.1113
$ws.Range("M68").Value = -4397.6875
$ws.Range("N68").Value = -7229.1113
$ws.Range("H71").Value = 5357.08
$ws.Range("I71").Value = 5146.6875
$ws.Range("J71").Value = 5731.1113
$ws.Range("K71").Value = 25733.4375
$ws.Range("L71").Value = 28655.5565
$ws.Range("M71").Value = -21989.4375
$ws.Range("N71").Value = -36143.5565
$ws.Range("H122").Value = 4134.76
$ws.Range("I122").Value = 2798.4443
$ws.Range("J122").Value = 7571
$ws.Range("K122").Value = 8395.332900000001
$ws.Range("L122").Value = 22713
$ws.Range("M122").Value = -5945.332900000001
$ws.Range("N122").Value = -27613
$ws.Range("H132").Value = 5503.2
$ws.Range("I132").Value = 3271.7144
$ws.Range("K132").Value = 9815.143199999999
$ws.Range("M132").Value = -7285.143199999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1998.5
$ws.Range("I126").Value = 1878.7
$ws.Range("J126").Value = 2597.5
$ws.Range("K126").Value = 5636.1
$ws.Range("L126").Value = 7792.5
$ws.Range("M126").Value = -3166.1
$ws.Range("N126").Value = -12732.5
